# Automated map update: remove the case "6269" record (row 65) from the
# PEBCOM sheet. All subsequent rows shift up by one to fill the gap,
# which matches the canonical diff (row 70 disappears, dimension shrinks
# from A1:P70 to A1:P69).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Delete the entire row 65; Excel automatically shifts rows 66:70 up to 65:69.
$ws.Rows.Item(65).Delete()
